# "Adapt tests to control version" - add a "version" column to the
# settings sheet of the form, with a constant value of 1, so downstream
# tests can pin/control the form version.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")
$ws.Activate()

# Existing header row: A1=form_title, B1=form_id -> add C1=version
$ws.Range("C1").Value = "version"

# Existing data row: A2=<title>, B2=<form id> -> add C2=1 (numeric)
$ws.Range("C2").Value = 1

# Move the active selection to the new last-used cell, like a user who
# just finished typing the new column would leave it.
$ws.Range("C3").Select() | Out-Null
